# Append two new ticket rows (143, 144) to the tickets log sheet.
# The date/time-looking values must stay literal text (matching the
# existing rows in the sheet), so force a Text number format on the
# date/time columns before writing them - otherwise Excel's
# autodetection would turn "2024-05-21" / "12:35:19" into date/time
# serial numbers instead of leaving them as strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A143:B144").NumberFormat = "@"
$ws.Range("H143:I144").NumberFormat = "@"

# Row 143
$ws.Range("A143").Value = "2024-05-21"
$ws.Range("B143").Value = "12:35:19"
$ws.Range("C143").Value = "Fallo en elevador"
$ws.Range("D143").Value = "-"
$ws.Range("E143").Value = "-"
$ws.Range("F143").Value = "-"
$ws.Range("G143").Value = "-"
$ws.Range("H143").Value = "12:35:21"
$ws.Range("I143").Value = "0:00:02"

# Row 144
$ws.Range("A144").Value = "2024-05-21"
$ws.Range("B144").Value = "12:36:05"
$ws.Range("C144").Value = "Etiquetadora21212"
$ws.Range("D144").Value = "-"
$ws.Range("E144").Value = "-"
$ws.Range("F144").Value = "-"
$ws.Range("G144").Value = "-"
$ws.Range("H144").Value = "12:36:06"
$ws.Range("I144").Value = "0:00:01"
